# Add a "Save" column (H) to the s_vals sheet, matching the existing
# header formatting used by the other header cells (B1:G1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1
$ws.Range("H1").Value = "Save"

# Copy the formatting (border/bold/alignment) from the neighbouring G1
# header cell onto the new H1 header cell.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Data values for H2:H13 taken from the commit diff
$saveValues = @(1, 0, 0, 1, 0, 1, 0, 1, 0, 1, 0, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
